# Auto-generated edit script: Add data for 2025-12-29
# Increments column L (2025 year-to-date totals) across Citywide Totals,
# By Neighborhood, and 32 individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6544  # was 6534 (Aggravated Assault)
$ws.Range("L3").Value = 7070  # was 7049 (Aggravated Battery)
$ws.Range("L4").Value = 1765  # was 1757 (Criminal Sexual Assault)
$ws.Range("L5").Value = 419  # was 418 (Homicide)
$ws.Range("L6").Value = 5789  # was 5779 (Robbery)
$ws.Range("L7").Value = 21587  # was 21537 (Total)

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 501  # was 500 (Aggravated Battery)
$ws.Range("L4").Value = 100  # was 99 (Criminal Sexual Assault)
$ws.Range("L7").Value = 1422  # was 1420 (Total)

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L4").Value = 20  # was 19 (Criminal Sexual Assault)
$ws.Range("L7").Value = 474  # was 473 (Total)

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 342  # was 341 (Aggravated Battery)
$ws.Range("L7").Value = 967  # was 966 (Total)

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 69  # was 68 (Robbery)
$ws.Range("L7").Value = 301  # was 300 (Total)

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 294  # was 292 (Aggravated Battery)
$ws.Range("L4").Value = 50  # was 49 (Criminal Sexual Assault)
$ws.Range("L6").Value = 214  # was 213 (Robbery)
$ws.Range("L7").Value = 832  # was 828 (Total)

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 134  # was 133 (Aggravated Battery)
$ws.Range("L7").Value = 432  # was 431 (Total)

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 79  # was 78 (Robbery)
$ws.Range("L7").Value = 374  # was 373 (Total)

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 78  # was 77 (Armour Square)
$ws.Range("L6").Value = 173  # was 172 (Ashburn)
$ws.Range("L7").Value = 692  # was 689 (Auburn Gresham)
$ws.Range("L8").Value = 1422  # was 1420 (Austin)
$ws.Range("L12").Value = 51  # was 50 (Beverly)
$ws.Range("L19").Value = 596  # was 594 (Chatham)
$ws.Range("L20").Value = 542  # was 541 (Chicago Lawn)
$ws.Range("L22").Value = 73  # was 72 (Clearing)
$ws.Range("L23").Value = 227  # was 226 (Douglas)
$ws.Range("L24").Value = 68  # was 66 (Dunning)
$ws.Range("L27").Value = 187  # was 186 (Edgewater)
$ws.Range("L29").Value = 1206  # was 1203 (Englewood)
$ws.Range("L33").Value = 967  # was 966 (Garfield Park)
$ws.Range("L37").Value = 832  # was 828 (Grand Crossing)
$ws.Range("L42").Value = 680  # was 679 (Humboldt Park)
$ws.Range("L47").Value = 151  # was 150 (Kenwood)
$ws.Range("L48").Value = 280  # was 279 (Lake View)
$ws.Range("L54").Value = 462  # was 459 (Loop)
$ws.Range("L65").Value = 432  # was 431 (New City)
$ws.Range("L67").Value = 750  # was 749 (North Lawndale)
$ws.Range("L76").Value = 343  # was 341 (River North)
$ws.Range("L78").Value = 283  # was 282 (Rogers Park)
$ws.Range("L79").Value = 601  # was 599 (Roseland)
$ws.Range("L83").Value = 474  # was 473 (South Chicago)
$ws.Range("L84").Value = 207  # was 206 (South Deering)
$ws.Range("L85").Value = 1075  # was 1072 (South Shore)
$ws.Range("L88").Value = 227  # was 226 (United Center)
$ws.Range("L94").Value = 263  # was 261 (West Loop)
$ws.Range("L95").Value = 301  # was 300 (West Pullman)
$ws.Range("L96").Value = 237  # was 236 (West Ridge)
$ws.Range("L97").Value = 174  # was 172 (West Town)
$ws.Range("L99").Value = 374  # was 373 (Woodlawn)
$ws.Range("L101").Value = 21587  # was 21537 (Total)

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 213  # was 212 (Aggravated Assault)
$ws.Range("L7").Value = 750  # was 749 (Total)

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 69  # was 68 (Aggravated Battery)
$ws.Range("L7").Value = 207  # was 206 (Total)

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 116  # was 114 (Aggravated Battery)
$ws.Range("L4").Value = 38  # was 37 (Criminal Sexual Assault)
$ws.Range("L7").Value = 462  # was 459 (Total)

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 466  # was 464 (Aggravated Battery)
$ws.Range("L6").Value = 288  # was 287 (Robbery)
$ws.Range("L7").Value = 1206  # was 1203 (Total)

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 54  # was 53 (Criminal Sexual Assault)
$ws.Range("L7").Value = 280  # was 279 (Total)

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L4").Value = 30  # was 29 (Criminal Sexual Assault)
$ws.Range("L6").Value = 165  # was 164 (Robbery)
$ws.Range("L7").Value = 596  # was 594 (Total)

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L5").Value = 10  # was 9 (Homicide)
$ws.Range("L6").Value = 154  # was 153 (Robbery)
$ws.Range("L7").Value = 343  # was 341 (Total)

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 70  # was 69 (Aggravated Assault)
$ws.Range("L7").Value = 173  # was 172 (Total)

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 184  # was 183 (Aggravated Assault)
$ws.Range("L7").Value = 680  # was 679 (Total)

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 74  # was 73 (Aggravated Assault)
$ws.Range("L7").Value = 283  # was 282 (Total)

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L6").Value = 16  # was 14 (Robbery)
$ws.Range("L7").Value = 68  # was 66 (Total)

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L4").Value = 21  # was 20 (Criminal Sexual Assault)
$ws.Range("L7").Value = 227  # was 226 (Total)

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 68  # was 67 (Aggravated Battery)
$ws.Range("L7").Value = 237  # was 236 (Total)

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 187  # was 185 (Aggravated Assault)
$ws.Range("L7").Value = 601  # was 599 (Total)

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 123  # was 122 (Robbery)
$ws.Range("L7").Value = 542  # was 541 (Total)

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 242  # was 240 (Aggravated Assault)
$ws.Range("L3").Value = 222  # was 221 (Aggravated Battery)
$ws.Range("L7").Value = 692  # was 689 (Total)

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 64  # was 62 (Aggravated Battery)
$ws.Range("L7").Value = 263  # was 261 (Total)

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 52  # was 51 (Aggravated Battery)
$ws.Range("L7").Value = 151  # was 150 (Total)

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 39  # was 37 (Aggravated Battery)
$ws.Range("L7").Value = 174  # was 172 (Total)

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L3").Value = 77  # was 76 (Aggravated Battery)
$ws.Range("L7").Value = 227  # was 226 (Total)

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 23  # was 22 (Aggravated Battery)
$ws.Range("L7").Value = 78  # was 77 (Total)

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L4").Value = 24  # was 23 (Criminal Sexual Assault)
$ws.Range("L7").Value = 187  # was 186 (Total)

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 320  # was 319 (Aggravated Assault)
$ws.Range("L3").Value = 445  # was 444 (Aggravated Battery)
$ws.Range("L6").Value = 226  # was 225 (Robbery)
$ws.Range("L7").Value = 1075  # was 1072 (Total)

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 24  # was 23 (Aggravated Assault)
$ws.Range("L7").Value = 73  # was 72 (Total)

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L3").Value = 16  # was 15 (Aggravated Battery)
$ws.Range("L7").Value = 51  # was 50 (Total)
